# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 3 (R) updates ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 230
$wsOff.Range("C3").Value = 154
$wsOff.Range("D3").Value = 43
$wsOff.Range("E3").Value = 20
$wsOff.Range("G3").Value = 3

# --- DEF sheet: row 3 (R) updates ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 233
$wsDef.Range("C3").Value = 173
$wsDef.Range("D3").Value = 49
$wsDef.Range("E3").Value = 23

$wb.Save()
